$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("A10:E10").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Посыл на Любовь"
$ws.Range("C11").Value = 43120
$ws.Range("D11").Value = "Посыл на Любовь подхватят везде!^В Посыле не бывает первых!^В Посыле и нации равны все,^Посыл – это промысел Верных!"
$ws.Range("E11").Value = "https://blagayavest.info/poems/20.01.18.html"
$ws.Hyperlinks.Add($ws.Range("E11"), "https://blagayavest.info/poems/20.01.18.html")
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)

# Row 12
$ws.Range("A11:E11").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Небо и небыль"
$ws.Range("C12").Value = 43256
$ws.Range("D12").Value = "Россия посылает в Небо Любовь!^Любовь отражается на континентах,^Посыл структурирует кровь,^Даже у своих оппонентов!"
$ws.Range("E12").Value = "https://blagayavest.info/poems/05.06.18.html"
$ws.Hyperlinks.Add($ws.Range("E12"), "https://blagayavest.info/poems/05.06.18.html")
$ws.Range("E11").Copy()
$ws.Range("E12").PasteSpecial(-4122)

# Row 13
$ws.Range("A12:E12").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Условия"
$ws.Range("C13").Value = 43313
$ws.Range("D13").Value = "Главное, Посыл и Решимость!^Желание Богов творить!^Проявится тогда и Терпимость,^Люди-Боги начнут говорить!"
$ws.Range("E13").Value = "https://blagayavest.info/poems/01.08.18.html"
$ws.Hyperlinks.Add($ws.Range("E13"), "https://blagayavest.info/poems/01.08.18.html")
$ws.Range("E12").Copy()
$ws.Range("E13").PasteSpecial(-4122)

# Row 14
$ws.Range("A13:E13").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Посыл и Промысел"
$ws.Range("C14").Value = 43333
$ws.Range("D14").Value = "Не хватает совсем немного!^Посыл должен быть полноценным!^Вам надо скорректировать слог,^Значит, сконцентрировать и Цель!"
$ws.Range("E14").Value = "https://blagayavest.info/poems/2018-08-21"
$ws.Hyperlinks.Add($ws.Range("E14"), "https://blagayavest.info/poems/2018-08-21")
$ws.Range("E13").Copy()
$ws.Range("E14").PasteSpecial(-4122)

# Row 15
$ws.Range("A14:E14").Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Рычаг Любви"
$ws.Range("C15").Value = 43334
$ws.Range("D15").Value = "Вы в Посылах несёте Любовь,^Информация проникает в тело,^Её к каждой клетке разносит кровь,^Человек в Любви – человек Света!"
$ws.Range("E15").Value = "https://blagayavest.info/poems/2018-08-22"
$ws.Hyperlinks.Add($ws.Range("E15"), "https://blagayavest.info/poems/2018-08-22")
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)

# Row 16
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Тысячи спасут миллионы"
$ws.Range("C16").Value = 43335
$ws.Range("D16").Value = "Повторяю, ваш Посыл действует!^Люди меняются, не зная того!^Главное, люди в Посыле участвуют,^И их неожиданно много!"
$ws.Range("E16").Value = "https://blagayavest.info/poems/2018-08-23"
$ws.Hyperlinks.Add($ws.Range("E16"), "https://blagayavest.info/poems/2018-08-23")
$ws.Range("E15").Copy()
$ws.Range("E16").PasteSpecial(-4122)

# Row 17
$ws.Range("A16:E16").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Человек - Золотое сечение"
$ws.Range("C17").Value = 43338
$ws.Range("D17").Value = "Посыл формирует Пространство!^Без Дуальности, только в Любви!^Человек прекращает странствия,^Энергии выравниваются внутри!"
$ws.Range("E17").Value = "https://blagayavest.info/poems/2018-08-26"
$ws.Hyperlinks.Add($ws.Range("E17"), "https://blagayavest.info/poems/2018-08-26")
$ws.Range("E16").Copy()
$ws.Range("E17").PasteSpecial(-4122)

# Row 18
$ws.Range("A17:E17").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Небывалое в Мире оружие"
$ws.Range("C18").Value = 43361
$ws.Range("D18").Value = "Выверенные вами Посылы,^Небывалое в Мире “оружие”!^Они сказками раньше слыли,^Теперь очищают они окружение!"
$ws.Range("E18").Value = "https://blagayavest.info/poems/2018-09-18"
$ws.Hyperlinks.Add($ws.Range("E18"), "https://blagayavest.info/poems/2018-09-18")
$ws.Range("E17").Copy()
$ws.Range("E18").PasteSpecial(-4122)

# Row 19
$ws.Range("A18:E18").Copy()
$ws.Range("A19:E19").PasteSpecial(-4122)
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Эта Чаша вас не минует"
$ws.Range("C19").Value = 43362
$ws.Range("D19").Value = "Посылы ведут к накоплению!^Достигнете Критической массы,^Произойдёт Историческое явление!^Появление людей-Богов касты!"
$ws.Range("E19").Value = "https://blagayavest.info/poems/2018-09-19"
$ws.Hyperlinks.Add($ws.Range("E19"), "https://blagayavest.info/poems/2018-09-19")
$ws.Range("E18").Copy()
$ws.Range("E19").PasteSpecial(-4122)

# Row 20
$ws.Range("A19:E19").Copy()
$ws.Range("A20:E20").PasteSpecial(-4122)
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "Сфера Разума"
$ws.Range("C20").Value = 43386
$ws.Range("D20").Value = "Отнеситесь к Посылу с Верой,^Влияние Любви увеличивается,^Посыл на Любовь – есть мера,^Тогда Спасение достигается!"
$ws.Range("E20").Value = "https://blagayavest.info/poems/2018-10-13"
$ws.Hyperlinks.Add($ws.Range("E20"), "https://blagayavest.info/poems/2018-10-13")
$ws.Range("E19").Copy()
$ws.Range("E20").PasteSpecial(-4122)

# Row 21
$ws.Range("A20:E20").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "В Небо пускают не всех"
$ws.Range("C21").Value = 43388
$ws.Range("D21").Value = "Каждый день, да ещё два раза,^Частота Посыла выше намного быта,^Не получится Преображения сразу,^Но, жизнь теперь не отбытие!"
$ws.Range("E21").Value = "https://blagayavest.info/poems/2018-10-15"
$ws.Hyperlinks.Add($ws.Range("E21"), "https://blagayavest.info/poems/2018-10-15")
$ws.Range("E20").Copy()
$ws.Range("E21").PasteSpecial(-4122)

# Row 22
$ws.Range("A21:E21").Copy()
$ws.Range("A22:E22").PasteSpecial(-4122)
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "Истоки и Приказ"
$ws.Range("C22").Value = 43419
$ws.Range("D22").Value = "Внешний прогресс – ничто!^Для Бога важнее – внутренний!^Человек, как и Бог, может всё!^Если Посыл – не ветреный!"
$ws.Range("E22").Value = "https://blagayavest.info/poems/2018-11-15"
$ws.Hyperlinks.Add($ws.Range("E22"), "https://blagayavest.info/poems/2018-11-15")
$ws.Range("E21").Copy()
$ws.Range("E22").PasteSpecial(-4122)

# Column width update for D
$ws.Columns.Item(4).ColumnWidth = 126.6667

# Update selection to D26
$ws.Range("D26").Select()

Write-Output "done"